$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Narrow column A from 147 to 108 (stored OOXML width units)
$ws.Range("A1").EntireColumn.ColumnWidth = 107.1

# 2. Swap the contents of row 3 and row 5 (values, not formatting)
$ws.Range("A3").Value = "('PL70 3380 5531 5257 6545 0735 9929', 'Multi-layered')"
$ws.Range("B3").Value = "[2019.0, 1.75088, 15.06291, 8.44]"
$ws.Range("C3").Value = "[2019.0, 1.75088, 15.06291, 78.44]"
$ws.Range("D3").Value = 70
$ws.Range("E3").Value = "[70.0]"

$ws.Range("A5").Value = "('HU52 9898 4213 5176 1777 7165 5419', 'encompassing')"
$ws.Range("B5").Value = "[2012.0, 1.7833, 23.57441, 23.3]"
$ws.Range("C5").Value = "[2012.0, 0.7833, 23.57441, 23.3]"
$ws.Range("D5").Value = -1
$ws.Range("E5").Value = "[-1.0]"

# 3. Update the "missing in ..." rows: swap entries 6/7 and 8/9, and
#    shorten the absolute path to the relative path.
$ws.Range("A6").Value = "('IS48 8566 8524 3637 6575 2319 52', '') missing in test_data/package_2\MOCK_DATA.xlsx"
$ws.Range("A7").Value = "('FR98 2851 7558 90QW BNYS BCF3 S20', '') missing in test_data/package_2\MOCK_DATA.xlsx"
$ws.Range("A8").Value = "('FR98 284551 7558 90QW BNYS BCF3 S20', '') missing in test_data/package_1\MOCK_DATA.xlsx"
$ws.Range("A9").Value = "('IS48 8566 8524 3637 6575 2319 52', 'Additional comment89') missing in test_data/package_1\MOCK_DATA.xlsx"
